$p = $ppt.ActivePresentation
$p.Slides.Item(10).Delete()
Write-Output $p.Slides.Count
